$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 6 new "Tomato" rows below the existing data (A9:A14)
for ($r = 9; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = "Tomato"
}

$ws.Range("A14").Select()
